# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps for the 95ad84ca-... row across
# the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" on Overview, and the matching
# "Correspond Handoff Datetime" on de-de share the same text value, so both
# must be updated together to stay merged as the same shared string.
$wsOverview.Range("G4").Value = "2016-08-15 10:44:56"
$wsDeDe.Range("H4").Value     = "2016-08-15 10:44:56"

# zh-cn row for 95ad84ca-...: Correspond Handoff Datetime / Handback DateTime
$wsZhCn.Range("H4").Value = "2016-08-15 10:44:51"
$wsZhCn.Range("K4").Value = "2016-08-15 10:45:17"

# de-de row for 95ad84ca-...: Correspond Handback DateTime
$wsDeDe.Range("K4").Value = "2016-08-15 10:45:24"
